$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

# Update Sheet2 InitialBalance cell (C2) from number 10000 to text "a1000"
# (set before the Sheet1 change so the shared-string table order matches:
#  "a1000" reuses the old "firefox" slot, "chrome" becomes the new entry)
$ws2.Range("C2").Value = "a1000"

# Update Sheet1 BrowseName cell (A2) from "firefox" to "chrome"
$ws1.Range("A2").Value = "chrome"

# Update the selection/active-cell shown on each sheet
$ws2.Range("D8").Select()
$ws1.Select()
$ws1.Range("A5").Select()
